# Insert a new data row at position 143 (pushing existing rows 143-191 down to 144-192)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A143").EntireRow.Insert()

$ws.Range("A143").Value = 3
$ws.Range("B143").Value = "Femacal de La Calera"
$ws.Range("C143").Value = "Coquimbo"
$ws.Range("D143").Value = 44468
$ws.Range("E143").Value = 5
$ws.Range("F143").Value = 100112012
$ws.Range("G143").Value = "Espinaca"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 120
$ws.Range("K143").Value = 3000
$ws.Range("L143").Value = 3000
$ws.Range("M143").Value = 3000
$ws.Range("N143").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O143").Value = "Provincia de Quillota"
$ws.Range("P143").Value = 1000
$ws.Range("Q143").Value = 3
$ws.Range("R143").Value = "Hortaliza"
